$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 14, shifting existing rows 14-23 down to 15-24
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new weekly price record
$ws.Range("A14").Value = 1
$ws.Range("B14").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C14").Value = "Arica y Parinacota"
$ws.Range("D14").Value = 44944
$ws.Range("E14").Value = 15
$ws.Range("F14").Value = "Fruta"
$ws.Range("G14").Value = 100103
$ws.Range("H14").Value = "Frutos de hueso (carozo)"
$ws.Range("I14").Value = 100103002
$ws.Range("J14").Value = "Ciruela"
$ws.Range("K14").Value = "Larry Ann"
$ws.Range("L14").Value = "Segunda"
$ws.Range("M14").Value = 300
$ws.Range("N14").Value = 7000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 7500
$ws.Range("Q14").Value = "$/bandeja 18 kilos granel"
$ws.Range("R14").Value = "Región de O'Higgins"
$ws.Range("S14").Value = 417
$ws.Range("T14").Value = 18
